$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.242.97"
$ws.Range("E2").Value = "  -2.97%  "
$ws.Range("D3").Value = "1.549.86"
$ws.Range("E3").Value = "  -4.97%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "206.55"
$ws.Range("E5").Value = "  -3.60%  "
$ws.Range("D6").Value = "1.01"
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("E7").Value = "  -5.39%  "
$ws.Range("E8").Value = "  -1.87%  "
$ws.Range("E9").Value = "  -3.49%  "
$ws.Range("D10").Value = "17.62"
$ws.Range("E10").Value = "  -5.15%  "
$ws.Range("D11").Value = "0.0778"
$ws.Range("E11").Value = "  -1.39%  "
$ws.Range("D12").Value = "1.766.03"
$ws.Range("E12").Value = "  -4.90%  "
$ws.Range("D13").Value = "1.554.53"
$ws.Range("E13").Value = "  -4.77%  "
$ws.Range("E14").Value = "  -4.99%  "
$ws.Range("D15").Value = "0.502"
$ws.Range("E15").Value = "  -4.85%  "
$ws.Range("D16").Value = "25.244.00"
$ws.Range("E16").Value = "  -2.94%  "
$ws.Range("E17").Value = "  -5.01%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "58.40"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.08%  "
$ws.Range("E19").Value = "  -0.07%  "
$ws.Range("D20").Value = "184.98"
$ws.Range("E20").Value = "  -4.10%  "
$ws.Range("E21").Value = "  -3.92%  "
$ws.Range("E22").Value = "  -3.33%  "
$ws.Range("D23").Value = "5.82"
$ws.Range("E23").Value = "  -4.17%  "
$ws.Range("E24").Value = "  -4.40%  "
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("D26").Value = "139.17"
$ws.Range("E26").Value = "  -3.56%  "
$ws.Range("E27").Value = "  -5.67%  "
$ws.Range("E28").Value = "  -3.17%  "
$ws.Range("E29").Value = "  -5.43%  "
$ws.Range("E30").Value = "  -7.25%  "
$ws.Range("E31").Value = "  -4.46%  "
$ws.Range("E32").Value = "  -3.68%  "
$ws.Range("E33").Value = "  -5.46%  "
$ws.Range("E34").Value = "  -3.40%  "
$ws.Range("E35").Value = "  -3.54%  "
$ws.Range("D36").Value = "1.080.54"
$ws.Range("E36").Value = "  -3.62%  "
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("E38").Value = "  -3.47%  "
$ws.Range("E39").Value = "  -5.46%  "
$ws.Range("E40").Value = "  -7.94%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.760"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -10.84%  "
$ws.Range("D42").Value = "0.798"
$ws.Range("E42").Value = "  +3.60%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "92.50"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.95%  "
$ws.Range("E44").Value = "  -2.41%  "
$ws.Range("D45").Value = "1.680.38"
$ws.Range("E45").Value = "  -4.89%  "
$ws.Range("E46").Value = "  -2.75%  "
$ws.Range("E47").Value = "  -2.51%  "
$ws.Range("D48").Value = "52.14"
$ws.Range("E48").Value = "  -4.47%  "
$ws.Range("E49").Value = "  -5.07%  "
$ws.Range("E51").Value = "  -2.10%  "

Write-Host "Updated cryptos list"